$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the single worker record in row 16 with the data that used to be
# in row 17 (document number, name, period) - the database refresh.
$ws.Range("C16").Value2 = 45547300
$ws.Range("D16").Value2 = "YERLING LAURINA POSADA BUELVAS"
$ws.Range("E16").Value2 = 1809

# Remove the now-duplicate second worker row (old row 17); this shifts the
# signature block rows (22/23 -> 21/22) up automatically.
$ws.Rows("17").Delete()

# Update summary figures to reflect the single remaining record.
$ws.Range("E11").Value2 = 31249
$ws.Range("C13").Value2 = 1
$ws.Range("F13").Value2 = 1
